# operation input tables updated
# A new row is inserted into the SpaceCoolingTechnology table (after the
# header row), pushing the existing data row down by one. The new row
# carries ID_SpaceCoolingTechnology=1, efficiency=3, power=10000, power_unit="W".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (the lone data row), shifting it
# down to row 3.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new technology entry.
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 3
$ws.Cells.Item(2, 3).Value = 10000
$ws.Cells.Item(2, 4).Value = "W"
